$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.011.58"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").Value = "1.642.40"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  +0.73%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "216.26"
$r.Style = "Normal"
$ws.Range("E5").Value = "  -0.21%  "
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "0.507"
$r.Style = "Normal"
$ws.Range("E6").Value = "  +0.52%  "
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "1.01"
$r.Style = "Normal"
$ws.Range("E7").Value = "  +0.62%  "
$ws.Range("E8").Value = "  -0.41%  "
$ws.Range("E9").Value = "  +0.42%  "
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "19.61"
$r.Style = "Normal"
$ws.Range("E10").Value = "  -0.62%  "
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "0.0796"
$r.Style = "Normal"
$ws.Range("E11").Value = "  +0.65%  "
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "1.869.13"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "4.29"
$r.Style = "Normal"
$ws.Range("E13").Value = "  -0.05%  "
$ws.Range("D14").Value = "1.615.46"
$ws.Range("E14").Value = "  -1.83%  "
$ws.Range("E15").Value = "  +0.09%  "
$ws.Range("D16").Value = "0.0₃0766"
$ws.Range("E16").Value = "  +0.41%  "
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "62.99"
$r.Style = "Normal"
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("D18").Value = "25.920.27"
$ws.Range("E18").Value = "  -0.86%  "
$ws.Range("E19").Value = "  +0.64%  "
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "192.96"
$r.Style = "Normal"
$ws.Range("E20").Value = "  -0.96%  "
$ws.Range("E21").Value = "  -1.68%  "
$ws.Range("E22").Value = "  -1.02%  "
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "1.81"
$r.Style = "Normal"
$ws.Range("E24").Value = "  +1.51%  "
$ws.Range("B25").Value = "Stellar"
$ws.Range("C25").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "0.131"
$r.Style = "Normal"
$ws.Range("E25").Value = "  +5.46%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "144.73"
$r.Style = "Normal"
$ws.Range("E26").Value = "  +1.53%  "
$ws.Range("E27").Value = "  +0.62%  "
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "6.93"
$r.Style = "Normal"
$ws.Range("E28").Value = "  +0.01%  "
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "15.56"
$r.Style = "Normal"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("E30").Value = "  -0.03%  "
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "0.0500"
$r.Style = "Normal"
$ws.Range("E31").Value = "  -0.41%  "
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "3.29"
$r.Style = "Normal"
$ws.Range("E32").Value = "  -1.90%  "
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "3.25"
$r.Style = "Normal"
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("E34").Value = "  -3.32%  "
$ws.Range("E35").Value = "  +2.59%  "
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "0.903"
$r.Style = "Normal"
$ws.Range("E36").Value = "  -0.90%  "
$ws.Range("D37").Value = "1.133.25"
$ws.Range("E37").Value = "  +0.18%  "
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "0.544"
$r.Style = "Normal"
$ws.Range("E38").Value = "  -1.54%  "
$ws.Range("E39").Value = "  -1.73%  "
$ws.Range("E40").Value = "  +0.06%  "
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "5.51"
$r.Style = "Normal"
$ws.Range("E41").Value = "  +0.45%  "
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "99.60"
$r.Style = "Normal"
$ws.Range("E42").Value = "  -0.75%  "
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "0.799"
$r.Style = "Normal"
$ws.Range("E43").Value = "  +0.25%  "
$ws.Range("D44").Value = "1.777.94"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").Value = "0.0₆0115"
$ws.Range("E45").Value = "  +3.98%  "
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "56.73"
$r.Style = "Normal"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("E47").Value = "  +2.78%  "
$ws.Range("E48").Value = "  -1.54%  "
$ws.Range("E49").Value = "  +0.34%  "
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "0.416"
$r.Style = "Normal"
$ws.Range("E50").Value = "  -0.15%  "
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "0.0961"
$r.Style = "Normal"
$ws.Range("E51").Value = "  -0.52%  "
